$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (number formats/styles) from the last existing row (1168) down to the new rows (1169:1181)
$ws.Range("A1168:V1168").Copy() | Out-Null
$ws.Range("A1169:V1181").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Row 1169
$ws.Cells.Item(1169, 1).Value = 'Entrainement'
$ws.Cells.Item(1169, 2).Value = 46041
$ws.Cells.Item(1169, 3).Value = 'Global'
$ws.Cells.Item(1169, 4).Value = 'J+2'
$ws.Cells.Item(1169, 5).Value = 'Mattheo Haon'
$ws.Cells.Item(1169, 6).Value = 'right back'
$ws.Cells.Item(1169, 7).Value = '01:12:31'
$ws.Cells.Item(1169, 8).Value = 5.67
$ws.Cells.Item(1169, 9).Value = 0.12
$ws.Cells.Item(1169, 10).Value = 5.54
$ws.Cells.Item(1169, 11).Value = 0.1
$ws.Cells.Item(1169, 12).Value = 0.02
$ws.Cells.Item(1169, 13).Value = 0
$ws.Cells.Item(1169, 14).Value = 0
$ws.Cells.Item(1169, 15).Value = 0
$ws.Cells.Item(1169, 16).Value = 4.59
$ws.Cells.Item(1169, 17).Value = 22.91
$ws.Cells.Item(1169, 18).Value = 5.12
$ws.Cells.Item(1169, 19).Value = 15
$ws.Cells.Item(1169, 20).Value = 5
$ws.Cells.Item(1169, 21).Value = 11
$ws.Cells.Item(1169, 22).Value = 4

# Row 1170
$ws.Cells.Item(1170, 1).Value = 'Entrainement'
$ws.Cells.Item(1170, 2).Value = 46041
$ws.Cells.Item(1170, 3).Value = 'Global'
$ws.Cells.Item(1170, 4).Value = 'J+2'
$ws.Cells.Item(1170, 5).Value = 'Ilan Ihaddadene'
$ws.Cells.Item(1170, 6).Value = 'center midfield'
$ws.Cells.Item(1170, 7).Value = '01:11:20'
$ws.Cells.Item(1170, 8).Value = 6.18
$ws.Cells.Item(1170, 9).Value = 0.1
$ws.Cells.Item(1170, 10).Value = 6.07
$ws.Cells.Item(1170, 11).Value = 0.1
$ws.Cells.Item(1170, 12).Value = 0
$ws.Cells.Item(1170, 13).Value = 0
$ws.Cells.Item(1170, 14).Value = 0
$ws.Cells.Item(1170, 15).Value = 0
$ws.Cells.Item(1170, 16).Value = 5.09
$ws.Cells.Item(1170, 17).Value = 19.98
$ws.Cells.Item(1170, 18).Value = 3.63
$ws.Cells.Item(1170, 19).Value = 13
$ws.Cells.Item(1170, 20).Value = 0
$ws.Cells.Item(1170, 21).Value = 7
$ws.Cells.Item(1170, 22).Value = 0

# Row 1171
$ws.Cells.Item(1171, 1).Value = 'Entrainement'
$ws.Cells.Item(1171, 2).Value = 46042
$ws.Cells.Item(1171, 3).Value = 'Global'
$ws.Cells.Item(1171, 4).Value = 'J+3'
$ws.Cells.Item(1171, 5).Value = 'Kamal Bafounta'
$ws.Cells.Item(1171, 6).Value = 'center midfield'
$ws.Cells.Item(1171, 7).Value = '01:44:52'
$ws.Cells.Item(1171, 8).Value = 8.2100000000000009
$ws.Cells.Item(1171, 9).Value = 1.44
$ws.Cells.Item(1171, 10).Value = 6.76
$ws.Cells.Item(1171, 11).Value = 0.78
$ws.Cells.Item(1171, 12).Value = 0.67
$ws.Cells.Item(1171, 13).Value = 0
$ws.Cells.Item(1171, 14).Value = 0
$ws.Cells.Item(1171, 15).Value = 0
$ws.Cells.Item(1171, 16).Value = 4.62
$ws.Cells.Item(1171, 17).Value = 22.75
$ws.Cells.Item(1171, 18).Value = 4.58
$ws.Cells.Item(1171, 19).Value = 49
$ws.Cells.Item(1171, 20).Value = 12
$ws.Cells.Item(1171, 21).Value = 31
$ws.Cells.Item(1171, 22).Value = 5

# Row 1172
$ws.Cells.Item(1172, 1).Value = 'Entrainement'
$ws.Cells.Item(1172, 2).Value = 46042
$ws.Cells.Item(1172, 3).Value = 'Global'
$ws.Cells.Item(1172, 4).Value = 'J+3'
$ws.Cells.Item(1172, 5).Value = 'Hedi Nasri'
$ws.Cells.Item(1172, 6).Value = 'right back'
$ws.Cells.Item(1172, 7).Value = '01:45:57'
$ws.Cells.Item(1172, 8).Value = 10.49
$ws.Cells.Item(1172, 9).Value = 1.31
$ws.Cells.Item(1172, 10).Value = 9.17
$ws.Cells.Item(1172, 11).Value = 1.18
$ws.Cells.Item(1172, 12).Value = 0.14000000000000001
$ws.Cells.Item(1172, 13).Value = 0
$ws.Cells.Item(1172, 14).Value = 0
$ws.Cells.Item(1172, 15).Value = 0
$ws.Cells.Item(1172, 16).Value = 4.3600000000000003
$ws.Cells.Item(1172, 17).Value = 24.21
$ws.Cells.Item(1172, 18).Value = 4.5599999999999996
$ws.Cells.Item(1172, 19).Value = 39
$ws.Cells.Item(1172, 20).Value = 7
$ws.Cells.Item(1172, 21).Value = 25
$ws.Cells.Item(1172, 22).Value = 9

# Row 1173
$ws.Cells.Item(1173, 1).Value = 'Entrainement'
$ws.Cells.Item(1173, 2).Value = 46042
$ws.Cells.Item(1173, 3).Value = 'Global'
$ws.Cells.Item(1173, 4).Value = 'J+3'
$ws.Cells.Item(1173, 5).Value = 'Karahali Souaré'
$ws.Cells.Item(1173, 6).Value = 'right forward'
$ws.Cells.Item(1173, 7).Value = '01:44:44'
$ws.Cells.Item(1173, 8).Value = 8.09
$ws.Cells.Item(1173, 9).Value = 1.34
$ws.Cells.Item(1173, 10).Value = 6.73
$ws.Cells.Item(1173, 11).Value = 1.1299999999999999
$ws.Cells.Item(1173, 12).Value = 0.24
$ws.Cells.Item(1173, 13).Value = 0
$ws.Cells.Item(1173, 14).Value = 0
$ws.Cells.Item(1173, 15).Value = 0
$ws.Cells.Item(1173, 16).Value = 4.28
$ws.Cells.Item(1173, 17).Value = 24.61
$ws.Cells.Item(1173, 18).Value = 5.76
$ws.Cells.Item(1173, 19).Value = 108
$ws.Cells.Item(1173, 20).Value = 33
$ws.Cells.Item(1173, 21).Value = 67
$ws.Cells.Item(1173, 22).Value = 17

# Row 1174
$ws.Cells.Item(1174, 1).Value = 'Entrainement'
$ws.Cells.Item(1174, 2).Value = 46042
$ws.Cells.Item(1174, 3).Value = 'Global'
$ws.Cells.Item(1174, 4).Value = 'J+3'
$ws.Cells.Item(1174, 5).Value = 'Romain Thunet'
$ws.Cells.Item(1174, 6).Value = 'center back'
$ws.Cells.Item(1174, 7).Value = '01:46:43'
$ws.Cells.Item(1174, 8).Value = 8.6
$ws.Cells.Item(1174, 9).Value = 0.16
$ws.Cells.Item(1174, 10).Value = 8.44
$ws.Cells.Item(1174, 11).Value = 0.15
$ws.Cells.Item(1174, 12).Value = 0.01
$ws.Cells.Item(1174, 13).Value = 0
$ws.Cells.Item(1174, 14).Value = 0
$ws.Cells.Item(1174, 15).Value = 0
$ws.Cells.Item(1174, 16).Value = 4.2
$ws.Cells.Item(1174, 17).Value = 21.95
$ws.Cells.Item(1174, 18).Value = 4.7300000000000004
$ws.Cells.Item(1174, 19).Value = 49
$ws.Cells.Item(1174, 20).Value = 2
$ws.Cells.Item(1174, 21).Value = 33
$ws.Cells.Item(1174, 22).Value = 7

# Row 1175
$ws.Cells.Item(1175, 1).Value = 'Entrainement'
$ws.Cells.Item(1175, 2).Value = 46042
$ws.Cells.Item(1175, 3).Value = 'Global'
$ws.Cells.Item(1175, 4).Value = 'J+3'
$ws.Cells.Item(1175, 5).Value = 'Theo Owono'
$ws.Cells.Item(1175, 6).Value = 'center midfield'
$ws.Cells.Item(1175, 7).Value = '01:46:25'
$ws.Cells.Item(1175, 8).Value = 7.48
$ws.Cells.Item(1175, 9).Value = 0.93
$ws.Cells.Item(1175, 10).Value = 6.55
$ws.Cells.Item(1175, 11).Value = 0.92
$ws.Cells.Item(1175, 12).Value = 0.01
$ws.Cells.Item(1175, 13).Value = 0
$ws.Cells.Item(1175, 14).Value = 0
$ws.Cells.Item(1175, 15).Value = 0
$ws.Cells.Item(1175, 16).Value = 4.17
$ws.Cells.Item(1175, 17).Value = 20.5
$ws.Cells.Item(1175, 18).Value = 4.09
$ws.Cells.Item(1175, 19).Value = 44
$ws.Cells.Item(1175, 20).Value = 2
$ws.Cells.Item(1175, 21).Value = 24
$ws.Cells.Item(1175, 22).Value = 4

# Row 1176
$ws.Cells.Item(1176, 1).Value = 'Entrainement'
$ws.Cells.Item(1176, 2).Value = 46042
$ws.Cells.Item(1176, 3).Value = 'Global'
$ws.Cells.Item(1176, 4).Value = 'J+3'
$ws.Cells.Item(1176, 5).Value = 'Malik Boussaid'
$ws.Cells.Item(1176, 6).Value = 'right back'
$ws.Cells.Item(1176, 7).Value = '01:44:44'
$ws.Cells.Item(1176, 8).Value = 8.0299999999999994
$ws.Cells.Item(1176, 9).Value = 1.26
$ws.Cells.Item(1176, 10).Value = 6.77
$ws.Cells.Item(1176, 11).Value = 0.79
$ws.Cells.Item(1176, 12).Value = 0.47
$ws.Cells.Item(1176, 13).Value = 0
$ws.Cells.Item(1176, 14).Value = 0
$ws.Cells.Item(1176, 15).Value = 0
$ws.Cells.Item(1176, 16).Value = 4.07
$ws.Cells.Item(1176, 17).Value = 24.74
$ws.Cells.Item(1176, 18).Value = 5.32
$ws.Cells.Item(1176, 19).Value = 32
$ws.Cells.Item(1176, 20).Value = 11
$ws.Cells.Item(1176, 21).Value = 32
$ws.Cells.Item(1176, 22).Value = 3

# Row 1177
$ws.Cells.Item(1177, 1).Value = 'Entrainement'
$ws.Cells.Item(1177, 2).Value = 46042
$ws.Cells.Item(1177, 3).Value = 'Global'
$ws.Cells.Item(1177, 4).Value = 'J+3'
$ws.Cells.Item(1177, 5).Value = 'Naim Ighbane'
$ws.Cells.Item(1177, 6).Value = 'center back'
$ws.Cells.Item(1177, 7).Value = '01:43:40'
$ws.Cells.Item(1177, 8).Value = 8.59
$ws.Cells.Item(1177, 9).Value = 1.3
$ws.Cells.Item(1177, 10).Value = 7.27
$ws.Cells.Item(1177, 11).Value = 0.85
$ws.Cells.Item(1177, 12).Value = 0.46
$ws.Cells.Item(1177, 13).Value = 0
$ws.Cells.Item(1177, 14).Value = 0
$ws.Cells.Item(1177, 15).Value = 1
$ws.Cells.Item(1177, 16).Value = 4.09
$ws.Cells.Item(1177, 17).Value = 25.22
$ws.Cells.Item(1177, 18).Value = 5.14
$ws.Cells.Item(1177, 19).Value = 55
$ws.Cells.Item(1177, 20).Value = 4
$ws.Cells.Item(1177, 21).Value = 31
$ws.Cells.Item(1177, 22).Value = 5

# Row 1178
$ws.Cells.Item(1178, 1).Value = 'Entrainement'
$ws.Cells.Item(1178, 2).Value = 46042
$ws.Cells.Item(1178, 3).Value = 'Global'
$ws.Cells.Item(1178, 4).Value = 'J+3'
$ws.Cells.Item(1178, 5).Value = 'Sofiane Belle'
$ws.Cells.Item(1178, 6).Value = 'left forward'
$ws.Cells.Item(1178, 7).Value = '01:46:32'
$ws.Cells.Item(1178, 8).Value = 7.83
$ws.Cells.Item(1178, 9).Value = 1.1299999999999999
$ws.Cells.Item(1178, 10).Value = 6.69
$ws.Cells.Item(1178, 11).Value = 1.1200000000000001
$ws.Cells.Item(1178, 12).Value = 0.02
$ws.Cells.Item(1178, 13).Value = 0
$ws.Cells.Item(1178, 14).Value = 0
$ws.Cells.Item(1178, 15).Value = 0
$ws.Cells.Item(1178, 16).Value = 3.53
$ws.Cells.Item(1178, 17).Value = 22.09
$ws.Cells.Item(1178, 18).Value = 5.1100000000000003
$ws.Cells.Item(1178, 19).Value = 26
$ws.Cells.Item(1178, 20).Value = 6
$ws.Cells.Item(1178, 21).Value = 20
$ws.Cells.Item(1178, 22).Value = 7

# Row 1179
$ws.Cells.Item(1179, 1).Value = 'Entrainement'
$ws.Cells.Item(1179, 2).Value = 46042
$ws.Cells.Item(1179, 3).Value = 'Global'
$ws.Cells.Item(1179, 4).Value = 'J+3'
$ws.Cells.Item(1179, 5).Value = 'Mattheo Haon'
$ws.Cells.Item(1179, 6).Value = 'right back'
$ws.Cells.Item(1179, 7).Value = '01:45:58'
$ws.Cells.Item(1179, 8).Value = 8.25
$ws.Cells.Item(1179, 9).Value = 1.21
$ws.Cells.Item(1179, 10).Value = 7.03
$ws.Cells.Item(1179, 11).Value = 1.07
$ws.Cells.Item(1179, 12).Value = 0.15
$ws.Cells.Item(1179, 13).Value = 0
$ws.Cells.Item(1179, 14).Value = 0
$ws.Cells.Item(1179, 15).Value = 0
$ws.Cells.Item(1179, 16).Value = 4.57
$ws.Cells.Item(1179, 17).Value = 22.41
$ws.Cells.Item(1179, 18).Value = 4.51
$ws.Cells.Item(1179, 19).Value = 41
$ws.Cells.Item(1179, 20).Value = 5
$ws.Cells.Item(1179, 21).Value = 22
$ws.Cells.Item(1179, 22).Value = 2

# Row 1180
$ws.Cells.Item(1180, 1).Value = 'Entrainement'
$ws.Cells.Item(1180, 2).Value = 46042
$ws.Cells.Item(1180, 3).Value = 'Global'
$ws.Cells.Item(1180, 4).Value = 'J+3'
$ws.Cells.Item(1180, 5).Value = 'Omar Benyounes'
$ws.Cells.Item(1180, 6).Value = 'center midfield'
$ws.Cells.Item(1180, 7).Value = '01:43:58'
$ws.Cells.Item(1180, 8).Value = 9.43
$ws.Cells.Item(1180, 9).Value = 2.44
$ws.Cells.Item(1180, 10).Value = 6.96
$ws.Cells.Item(1180, 11).Value = 1.95
$ws.Cells.Item(1180, 12).Value = 0.52
$ws.Cells.Item(1180, 13).Value = 0
$ws.Cells.Item(1180, 14).Value = 0
$ws.Cells.Item(1180, 15).Value = 0
$ws.Cells.Item(1180, 16).Value = 4.92
$ws.Cells.Item(1180, 17).Value = 24.46
$ws.Cells.Item(1180, 18).Value = 5.22
$ws.Cells.Item(1180, 19).Value = 81
$ws.Cells.Item(1180, 20).Value = 21
$ws.Cells.Item(1180, 21).Value = 43
$ws.Cells.Item(1180, 22).Value = 17

# Row 1181
$ws.Cells.Item(1181, 1).Value = 'Entrainement'
$ws.Cells.Item(1181, 2).Value = 46042
$ws.Cells.Item(1181, 3).Value = 'Global'
$ws.Cells.Item(1181, 4).Value = 'J+3'
$ws.Cells.Item(1181, 5).Value = 'Ilan Ihaddadene'
$ws.Cells.Item(1181, 6).Value = 'center midfield'
$ws.Cells.Item(1181, 7).Value = '01:44:35'
$ws.Cells.Item(1181, 8).Value = 8.36
$ws.Cells.Item(1181, 9).Value = 1.29
$ws.Cells.Item(1181, 10).Value = 7.07
$ws.Cells.Item(1181, 11).Value = 1.28
$ws.Cells.Item(1181, 12).Value = 0.01
$ws.Cells.Item(1181, 13).Value = 0
$ws.Cells.Item(1181, 14).Value = 0
$ws.Cells.Item(1181, 15).Value = 0
$ws.Cells.Item(1181, 16).Value = 4.7300000000000004
$ws.Cells.Item(1181, 17).Value = 20.81
$ws.Cells.Item(1181, 18).Value = 4.45
$ws.Cells.Item(1181, 19).Value = 39
$ws.Cells.Item(1181, 20).Value = 4
$ws.Cells.Item(1181, 21).Value = 16
$ws.Cells.Item(1181, 22).Value = 3

$ws.Range("C1184").Select() | Out-Null